# Add a new "DirectLine" column to the company resources sheet, between
# "Role" (F) and "MobilePhone1" (old G, now shifted to H).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at G; everything from the old G onward shifts
# right by one (MobilePhone1/2/3, WorkEmailAddress, PersonalEmailAddress,
# BusinessUnit, BudgetMonth all move one column over).
$ws.Columns("G:G").Insert()

# Header for the new column.
$ws.Range("G1").Value = "DirectLine"

# The freshly inserted G2:G5 cells picked up formatting from column F on
# insert; match them to the (empty, unfilled) style that the rest of that
# data block uses - the same style now sitting in H2:H5.
$ws.Range("H2:H5").Copy()
$ws.Range("G2:G5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Give the new column the same width as the Role column to its left.
$ws.Columns("G:G").ColumnWidth = $ws.Columns("F:F").ColumnWidth

# Restore the selection to the (now shifted) MobilePhone2 cell.
$ws.Range("I9").Select() | Out-Null
